$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DashboardWireFrame")
$r = $ws.Range("D13")
$r.NumberFormat = "mm-dd-yy"
$r.Font.Bold = $true
